$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District names (G column) to official names from the website
$ws.Range("G4").Value = "Vijayapura (Bijapur)"
$ws.Range("G5").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G14").Value = "Vijayapura (Bijapur)"
$ws.Range("G21").Value = "Vijayapura (Bijapur)"
$ws.Range("G24").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G27").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G29").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G30").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G32").Value = "Ballari (Bellary)"
$ws.Range("G34").Value = "Vijayapura (Bijapur)"
$ws.Range("G35").Value = "Chikkamagaluru (Chikmagalur)"

# Remove the stray empty Address cells that had no content
$ws.Range("F6").ClearContents()
$ws.Range("F19").ClearContents()
